$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 55359.223852715098
$ws.Range("C2").Value = 33694.234850465
$ws.Range("D2").Value = 2193.49549969895

$ws.Range("B3").Value = 69796.677594866298
$ws.Range("C3").Value = 46783.842901226097
$ws.Range("D3").Value = 845.64980830880904

$ws.Range("B4").Value = 80743.654556177498
$ws.Range("C4").Value = 57316.646779459297
$ws.Range("D4").Value = 431.47672523089102

$ws.Range("B5").Value = 92971.7131604608
$ws.Range("C5").Value = 69303.094813458898

$ws.Range("B6").Value = 100465.33189307401
$ws.Range("C6").Value = 76735.641829743996
$ws.Range("D6").Value = 128.79443861830501

$ws.Range("B4:D4").Select()
